$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 942
$ws1.Range("F5").Value = 225
$ws1.Range("F6").Value = 48
$ws1.Range("F7").Value = 1167
$ws1.Range("F8").Value = 934
$ws1.Range("F9").Value = 30
$ws1.Range("F10").Value = 725
$ws1.Range("F11").Value = 1041
$ws1.Range("F12").Value = 1484
$ws1.Range("F15").Value = 1653
$ws1.Range("F26").Value = 479
$ws1.Range("F32").Value = 2437
$ws1.Range("F34").Value = 1393
$ws1.Range("F35").Value = 466
$ws1.Range("F38").Value = 4003

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1039
$ws2.Range("F10").Value = 85
$ws2.Range("F36").Value = 2

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1277
$ws3.Range("F5").Value = 1674

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1277
$ws4.Range("F4").Value = 1674
$ws4.Range("F8").Value = 942
$ws4.Range("F10").Value = 48
$ws4.Range("F11").Value = 1167
$ws4.Range("F12").Value = 934
$ws4.Range("F14").Value = 30
$ws4.Range("F16").Value = 725
$ws4.Range("F20").Value = 1041
$ws4.Range("F21").Value = 1484
$ws4.Range("F24").Value = 1653
$ws4.Range("F34").Value = 479
$ws4.Range("F42").Value = 2437
$ws4.Range("F47").Value = 1393
$ws4.Range("F48").Value = 466
$ws4.Range("F50").Value = 4003
